$d = $word.ActiveDocument

$d.Content.Find.Execute("127×7=889", $true, $false, $false, $false, $false, $true, 1, $false, "516×7=3612", 2) | Out-Null
$d.Content.Find.Execute("791×4=3164", $true, $false, $false, $false, $false, $true, 1, $false, "425×2=850", 2) | Out-Null
$d.Content.Find.Execute("754×7=5278", $true, $false, $false, $false, $false, $true, 1, $false, "838×3=2514", 2) | Out-Null
$d.Content.Find.Execute("870×7=6090", $true, $false, $false, $false, $false, $true, 1, $false, "824×4=3296", 2) | Out-Null
$d.Content.Find.Execute("114×8=912", $true, $false, $false, $false, $false, $true, 1, $false, "531×7=3717", 2) | Out-Null
$d.Content.Find.Execute("150×3=450", $true, $false, $false, $false, $false, $true, 1, $false, "154×7=1078", 2) | Out-Null
$d.Content.Find.Execute("210×7=1470", $true, $false, $false, $false, $false, $true, 1, $false, "730×8=5840", 2) | Out-Null
$d.Content.Find.Execute("553×4=2212", $true, $false, $false, $false, $false, $true, 1, $false, "857×8=6856", 2) | Out-Null
$d.Content.Find.Execute("265×9=2385", $true, $false, $false, $false, $false, $true, 1, $false, "320×7=2240", 2) | Out-Null
$d.Content.Find.Execute("713×4=2852", $true, $false, $false, $false, $false, $true, 1, $false, "807×9=7263", 2) | Out-Null
$d.Content.Find.Execute("698×6=4188", $true, $false, $false, $false, $false, $true, 1, $false, "850×5=4250", 2) | Out-Null
$d.Content.Find.Execute("655×9=5895", $true, $false, $false, $false, $false, $true, 1, $false, "427×8=3416", 2) | Out-Null
$d.Content.Find.Execute("761×2=1522", $true, $false, $false, $false, $false, $true, 1, $false, "584×4=2336", 2) | Out-Null
$d.Content.Find.Execute("690×7=4830", $true, $false, $false, $false, $false, $true, 1, $false, "529×6=3174", 2) | Out-Null
$d.Content.Find.Execute("652×5=3260", $true, $false, $false, $false, $false, $true, 1, $false, "760×8=6080", 2) | Out-Null
$d.Content.Find.Execute("236×2=472", $true, $false, $false, $false, $false, $true, 1, $false, "935×9=8415", 2) | Out-Null
$d.Content.Find.Execute("911×2=1822", $true, $false, $false, $false, $false, $true, 1, $false, "416×6=2496", 2) | Out-Null
$d.Content.Find.Execute("746×3=2238", $true, $false, $false, $false, $false, $true, 1, $false, "568×6=3408", 2) | Out-Null
$d.Content.Find.Execute("392×9=3528", $true, $false, $false, $false, $false, $true, 1, $false, "731×3=2193", 2) | Out-Null
$d.Content.Find.Execute("684×7=4788", $true, $false, $false, $false, $false, $true, 1, $false, "259×7=1813", 2) | Out-Null
$d.Content.Find.Execute("696×7=4872", $true, $false, $false, $false, $false, $true, 1, $false, "599×8=4792", 2) | Out-Null
$d.Content.Find.Execute("485×8=3880", $true, $false, $false, $false, $false, $true, 1, $false, "241×2=482", 2) | Out-Null
$d.Content.Find.Execute("595×6=3570", $true, $false, $false, $false, $false, $true, 1, $false, "559×3=1677", 2) | Out-Null
$d.Content.Find.Execute("185×7=1295", $true, $false, $false, $false, $false, $true, 1, $false, "175×8=1400", 2) | Out-Null
$d.Content.Find.Execute("173×6=1038", $true, $false, $false, $false, $false, $true, 1, $false, "498×7=3486", 2) | Out-Null
